$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect Price/Volume columns as text so numeric-looking strings are not coerced to numbers
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.703.23'
$ws.Range("E2").Value = '  -1.32%  '

# Row 3
$ws.Range("D3").Value = '2.097.30'
$ws.Range("E3").Value = '  -0.39%  '

# Row 4
$ws.Range("E4").Value = '  +0.52%  '

# Row 5
$ws.Range("D5").Value = '343.53'
$ws.Range("E5").Value = '  -1.94%  '

# Row 6
$ws.Range("E6").Value = '  +0.60%  '

# Row 7
$ws.Range("D7").Value = '0.5160'
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").Value = '0.4380'
$ws.Range("E8").Value = '  -2.21%  '

# Row 9
$ws.Range("D9").Value = '53.45'
$ws.Range("E9").Value = '  +1.34%  '

# Row 10
$ws.Range("D10").Value = '0.09186'
$ws.Range("E10").Value = '  +2.59%  '

# Row 11
$ws.Range("D11").Value = '1.166'
$ws.Range("E11").Value = '  -0.97%  '

# Row 12
$ws.Range("E12").Value = '  -4.93%  '

# Row 13
$ws.Range("D13").Value = '2.088.84'
$ws.Range("E13").Value = '  -1.05%  '

# Row 14
$ws.Range("D14").Value = '6.764'
$ws.Range("E14").Value = '  -0.04%  '

# Row 15
$ws.Range("D15").Value = '8.175'
$ws.Range("E15").Value = '  -1.14%  '

# Row 16
$ws.Range("D16").Value = '102.51'
$ws.Range("E16").Value = '  +3.21%  '

# Row 17
$ws.Range("E17").Value = '  +0.28%  '

# Row 18
$ws.Range("E18").Value = '  +0.49%  '

# Row 19
$ws.Range("D19").Value = '21.02'
$ws.Range("E19").Value = '  +0.64%  '

# Row 20
$ws.Range("D20").Value = '0.06673'
$ws.Range("E20").Value = '  +0.18%  '

# Row 21
$ws.Range("D21").Value = '1.008'
$ws.Range("E21").Value = '  +0.56%  '

# Row 22
$ws.Range("D22").Value = '6.203'
$ws.Range("E22").Value = '  -1.18%  '

# Row 23
$ws.Range("D23").Value = '29.763.13'
$ws.Range("E23").Value = '  -1.45%  '

# Row 24
$ws.Range("D24").Value = '12.57'
$ws.Range("E24").Value = '  -2.92%  '

# Row 25
$ws.Range("D25").Value = '2.306'
$ws.Range("E25").Value = '  -2.09%  '

# Row 26
$ws.Range("D26").Value = '2.345.12'
$ws.Range("E26").Value = '  -0.57%  '

# Row 27
$ws.Range("D27").Value = '21.92'
$ws.Range("E27").Value = '  -0.59%  '

# Row 28
$ws.Range("D28").Value = '162.00'
$ws.Range("E28").Value = '  -0.73%  '

# Row 29
$ws.Range("D29").Value = '2.498'
$ws.Range("E29").Value = '  -2.22%  '

# Row 30
$ws.Range("D30").Value = '133.32'
$ws.Range("E30").Value = '  -0.33%  '

# Row 31
$ws.Range("D31").Value = '1.129'
$ws.Range("E31").Value = '  -4.67%  '

# Row 32
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '1.667'
$ws.Range("E32").Value = '  +0.97%  '

# Row 33
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.1051'
$ws.Range("E33").Value = '  -1.69%  '

# Row 34
$ws.Range("D34").Value = '6.193'
$ws.Range("E34").Value = '  -1.23%  '

# Row 35
$ws.Range("D35").Value = '3.963'
$ws.Range("E35").Value = '  -0.41%  '

# Row 36
$ws.Range("D36").Value = '6.306'
$ws.Range("E36").Value = '  +6.81%  '

# Row 37
$ws.Range("D37").Value = '10.39'
$ws.Range("E37").Value = '  +1.82%  '

# Row 38
$ws.Range("D38").Value = '0.02575'
$ws.Range("E38").Value = '  -0.74%  '

# Row 39
$ws.Range("D39").Value = '0.06692'
$ws.Range("E39").Value = '  -2.33%  '

# Row 40
$ws.Range("D40").Value = '0.6991'
$ws.Range("E40").Value = '  +1.91%  '

# Row 41
$ws.Range("D41").Value = '12.45'
$ws.Range("E41").Value = '  -0.96%  '

# Row 42
$ws.Range("D42").Value = '1.329'
$ws.Range("E42").Value = '  +6.06%  '

# Row 43
$ws.Range("D43").Value = '0.2216'
$ws.Range("E43").Value = '  -4.75%  '

# Row 44
$ws.Range("D44").Value = '0.6805'
$ws.Range("E44").Value = '  +5.51%  '

# Row 45
$ws.Range("D45").Value = '14.33'
$ws.Range("E45").Value = '  -0.12%  '

# Row 46
$ws.Range("D46").Value = '2.315'
$ws.Range("E46").Value = '  +0.49%  '

# Row 47
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000358'
$ws.Range("E47").Value = '  -3.08%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.625'
$ws.Range("E48").Value = '  -1.15%  '

# Row 49
$ws.Range("B49").Value = 'WEMIXTOKEN'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '1.209'
$ws.Range("E49").Value = '  +3.60%  '

# Row 50
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.217'
$ws.Range("E50").Value = '  -0.65%  '

# Row 51
$ws.Range("D51").Value = '81.20'
$ws.Range("E51").Value = '  -3.47%  '
